# Update "想去人数" (F column) figures on the "展览" sheet and the
# aggregated "全部类型" sheet, matching the refreshed gh-pages data dump.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsAll        = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 sheet (rows 2..26 correspond to individual listings) ---
$wsExhibition.Range("F3").Value  = 12836
$wsExhibition.Range("F5").Value  = 82
$wsExhibition.Range("F6").Value  = 64
$wsExhibition.Range("F7").Value  = 39
$wsExhibition.Range("F8").Value  = 19
$wsExhibition.Range("F9").Value  = 10
$wsExhibition.Range("F10").Value = 12750
$wsExhibition.Range("F11").Value = 275
$wsExhibition.Range("F12").Value = 23
$wsExhibition.Range("F13").Value = 8655
$wsExhibition.Range("F14").Value = 7645
$wsExhibition.Range("F15").Value = 188
$wsExhibition.Range("F17").Value = 430
$wsExhibition.Range("F18").Value = 123
$wsExhibition.Range("F19").Value = 977
$wsExhibition.Range("F20").Value = 6
$wsExhibition.Range("F24").Value = 14

# --- 全部类型 sheet (same listings, shifted down by one row because it
#     also includes the single "演出" entry at the top) ---
$wsAll.Range("F4").Value  = 12836
$wsAll.Range("F6").Value  = 82
$wsAll.Range("F7").Value  = 64
$wsAll.Range("F8").Value  = 39
$wsAll.Range("F9").Value  = 19
$wsAll.Range("F10").Value = 10
$wsAll.Range("F11").Value = 12750
$wsAll.Range("F12").Value = 275
$wsAll.Range("F13").Value = 23
$wsAll.Range("F14").Value = 8654
$wsAll.Range("F15").Value = 7645
$wsAll.Range("F16").Value = 188
$wsAll.Range("F18").Value = 430
$wsAll.Range("F19").Value = 123
$wsAll.Range("F20").Value = 977
$wsAll.Range("F21").Value = 6
$wsAll.Range("F26").Value = 14
